# docs: Added obsFile and moved rr subpackages in functionalities list for 0.3.0 release
# Refs: #274.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Source table")

# --- Observation station file (new) row (Source table row 48) ---
# Was only marked Read/Write = "O" (possibly future support) with no
# version/module/class info. Now it is fully supported as of 0.3.0, backed
# by the new ObservationPointModel.
$src.Range("B48").Value = "X"
$src.Range("C48").Value = "X"
$src.Range("D48").Value = "0.3.0"
$src.Range("E48").Value = "hydrolib.core.io.obs.models"
$src.Range("F48").Value = "ObservationPointModel"

# --- Main sobek_3b.fnm row (Source table row 59) ---
# RainfallRunoffModel moved from hydrolib.core.io.fnm.models into
# hydrolib.core.io.rr.models as of 0.3.0.
$src.Range("E59").Value = "hydrolib.core.io.rr.models"
$src.Range("G59").Value = "Used to be in hydrolib.core.io.fnm.models before 0.3.0"

# --- Rainfall .bui file row (Source table row 60) ---
# BuiModel moved from hydrolib.core.io.bui.models into
# hydrolib.core.io.rr.meteo.models as of 0.3.0.
$src.Range("G60").Value = "Used to be in hydrolib.core.io.bui.models before 0.3.0"
$src.Range("E60").Value = "hydrolib.core.io.rr.meteo.models"

# --- Update the view/selection state to match the author's last position ---
# Touch "Source table" first, then re-activate "FM mkdocs table" last so it
# remains the active (tabSelected) sheet, matching activeTab="2".
$src.Activate()
$src.Range("B29").Select()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1

$fm = $wb.Worksheets.Item("FM mkdocs table")
$fm.Activate()
$fm.Range("A53").Select()
